$d = $word.ActiveDocument

# --- Edit 1: append trailing spaces + a red "(This is a change ...)" note
#     to the first paragraph, split across three runs exactly like the
#     target markup (as Word would when the red note was typed in pieces).
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$r.End = $p1.Range.End - 1   # exclude the paragraph mark
$r.Collapse(0)                # 0 = wdCollapseEnd
$r.InsertAfter("  ")
$r.Collapse(0)

function Insert-ColoredRun($range, $text) {
    $range.InsertAfter($text)
    $startPos = $range.End - $text.Length
    $endPos = $range.End
    $colorRange = $range.Document.Range($startPos, $endPos)
    $colorRange.Font.Color = 255   # 255 = wdColorRed (RGB 0xFF0000)
    $range.Collapse(0)             # 0 = wdCollapseEnd
}

$enDash = [char]0x2013
$part1 = "(This is a change " + $enDash + " Ve"
$part2 = "rsion for main branch"
$part3 = ")"

Insert-ColoredRun $r $part1
Insert-ColoredRun $r $part2
Insert-ColoredRun $r $part3

# --- Edit 2: remove the trailing "ank God almighty, we are free at last."
#     paragraph that followed the end of "The Raven".
$last = $d.Paragraphs.Last
$last.Range.Delete()
